$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the "NBOX Air1" listing now appears where "NBOX MARATHON" used to be
$ws.Range("A2").Value = "NBOX Air1 TWS On Ear True Wireless (TWS) 20 Hours Playback IPX5(Splash & Sweat Proof) Passive noise cancellation -Bluetooth Black"

# Row 3: ... and "NBOX MARATHON" now appears where "NBOX Air1" used to be
$ws.Range("A3").Value = "NBOX MARATHON Over Ear Bluetooth Neckband 20 Hours Playback IPX5(Splash & Sweat Proof) Passive noise cancellation -Bluetooth Silver"

# Row 4: now holds the "NBOX G2 PLATINUM" listing with a refreshed price
$ws.Range("A4").Value = "NBOX G2 PLATINUM SERISE Over Ear Bluetooth Neckband 6 Hours Playback IPX5(Splash & Sweat Proof) Passive noise cancellation -Bluetooth Gray"
$ws.Range("B4").Value = "Rs. 700"

# Row 5: new "boAt Airdopes 171" listing with its price
$ws.Range("A5").Value = "boAt Airdopes 171 On Ear Wireless With Mic Headphones/Earphones Black"
$ws.Range("B5").Value = "Rs. 1,499"
